{"js": "// Apply the \"Added many more features\" edit:\n// - Retitle the review (title appears twice: H1 heading + bold recap line)\n// - Rewrite the \"What we like\" bullet list (4 bullets)\n// - Rewrite one \"What we don't like\" bullet\n// - Rewrite the italic summary/description paragraph near the end\n\nconst replacements = [\n  {\n    from: \"Play Inferno Free - Classic Fruit Machine Review\",\n    to: \"Play Inferno Free - Review of The Classic Fruit Machine\",\n  },\n  {\n    from: \"Visually pleasing design with bright colors\",\n    to: \"Visually pleasing with bright colors and well-defined fruits\",\n  },\n  {\n    from: \"Classic structure, perfect for traditional slot lovers\",\n    to: \"Classic fruit machine design\",\n  },\n  {\n    from: \"Demo version available to play for free\",\n    to: \"Autoplay mode available\",\n  },\n  {\n    from: \"'Leave or double' function adds excitement to gameplay\",\n    to: \"Demo version to play for free\",\n  },\n  {\n    from: \"Limited paylines\",\n    to: \"High volatility\",\n  },\n  {\n    from:\n      \"Read our review of Inferno, a classic fruit machine from Novomatic with 5 paylines and a high volatility. Play for free and enjoy the traditional design.\",\n    to:\n      \"Review of Inferno, a classic fruit machine with a fiery background. Play for free and enjoy the visually pleasing design.\",\n  },\n];\n\nfor (const { from, to } of replacements) {\n  const results = context.document.body.search(from, { matchCase: true });\n  results.load(\"text\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(to, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Apply the \"Added many more features\" edit:\n# - Retitle the review (title appears twice: H1 heading + bold recap line)\n# - Rewrite the \"What we like\" bullet list (4 bullets)\n# - Rewrite one \"What we don't like\" bullet\n# - Rewrite the italic summary/description paragraph near the end\n\n$d = $word.ActiveDocument\n\nfunction Replace-AllText($findText, $replaceText) {\n    $range = $d.Content\n    $find = $range.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $findText\n    $find.Replacement.Text = $replaceText\n    $find.Forward = $true\n    $find.Wrap = 1            # wdFindContinue\n    $find.MatchCase = $true\n    $find.MatchWholeWord = $false\n    $find.MatchWildcards = $false\n    $find.Execute([ref]$findText, $false, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null\n}\n\nReplace-AllText \"Play Inferno Free - Classic Fruit Machine Review\" \"Play Inferno Free - Review of The Classic Fruit Machine\"\nReplace-AllText \"Visually pleasing design with bright colors\" \"Visually pleasing with bright colors and well-defined fruits\"\nReplace-AllText \"Classic structure, perfect for traditional slot lovers\" \"Classic fruit machine design\"\nReplace-AllText \"Demo version available to play for free\" \"Autoplay mode available\"\nReplace-AllText \"'Leave or double' function adds excitement to gameplay\" \"Demo version to play for free\"\nReplace-AllText \"Limited paylines\" \"High volatility\"\nReplace-AllText \"Read our review of Inferno, a classic fruit machine from Novomatic with 5 paylines and a high volatility. Play for free and enjoy the traditional design.\" \"Review of Inferno, a classic fruit machine with a fiery background. Play for free and enjoy the visually pleasing design.\"\n"}
